$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.373.55"
$ws.Range("D3").Value = "'1.666.45"
$ws.Range("E3").Value = "'  +1.80%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  +0.18%  "
$ws.Range("D5").Value = "'311.79"
$ws.Range("E5").Value = "'  +1.92%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.07%  "
$ws.Range("D7").Value = "'0.3943"
$ws.Range("E7").Value = "'  +1.98%  "
$ws.Range("D8").Value = "'0.3914"
$ws.Range("E8").Value = "'  +2.40%  "
$ws.Range("E9").Value = "'  +5.99%  "
$ws.Range("D10").Value = "'1.385"
$ws.Range("E10").Value = "'  +3.42%  "
$ws.Range("D11").Value = "'1.001"
$ws.Range("E11").Value = "'  +0.10%  "
$ws.Range("D12").Value = "'0.08551"
$ws.Range("E12").Value = "'  +1.13%  "
$ws.Range("D13").Value = "'24.40"
$ws.Range("E13").Value = "'  +3.71%  "
$ws.Range("D14").Value = "'7.262"
$ws.Range("E14").Value = "'  +3.00%  "
$ws.Range("D15").Value = "'7.970"
$ws.Range("E15").Value = "'  +7.46%  "
$ws.Range("D16").Value = "'0.00001330"
$ws.Range("E16").Value = "'  +4.49%  "
$ws.Range("D17").Value = "'1.663.28"
$ws.Range("E17").Value = "'  +1.56%  "
$ws.Range("D18").Value = "'94.83"
$ws.Range("E18").Value = "'  -0.01%  "
$ws.Range("D19").Value = "'0.07022"
$ws.Range("E19").Value = "'  +2.27%  "
$ws.Range("D20").Value = "'20.53"
$ws.Range("E20").Value = "'  -0.31%  "
$ws.Range("D21").Value = "'6.983"
$ws.Range("E21").Value = "'  +1.77%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "'  +0.03%  "
$ws.Range("E23").Value = "'  +2.14%  "
$ws.Range("D24").Value = "'24.430.19"
$ws.Range("E24").Value = "'  +1.55%  "
$ws.Range("D25").Value = "'2.486"
$ws.Range("E25").Value = "'  +6.79%  "
$ws.Range("D26").Value = "'3.079"
$ws.Range("E26").Value = "'  +15.22%  "
$ws.Range("D27").Value = "'22.53"
$ws.Range("E27").Value = "'  +1.46%  "
$ws.Range("D28").Value = "'157.09"
$ws.Range("D29").Value = "'142.41"
$ws.Range("E29").Value = "'  +2.00%  "
$ws.Range("D30").Value = "'5.330"
$ws.Range("E30").Value = "'  +0.02%  "
$ws.Range("D31").Value = "'7.894"
$ws.Range("E31").Value = "'  -8.79%  "
$ws.Range("D32").Value = "'2.548"
$ws.Range("E32").Value = "'  +5.49%  "
$ws.Range("D33").Value = "'1.849.71"
$ws.Range("E33").Value = "'  +2.04%  "
$ws.Range("D34").Value = "'1.061"
$ws.Range("E34").Value = "'  +12.50%  "
$ws.Range("D35").Value = "'0.03109"
$ws.Range("E35").Value = "'  +8.40%  "
$ws.Range("D36").Value = "'0.08209"
$ws.Range("E36").Value = "'  +3.08%  "
$ws.Range("D37").Value = "'6.861"
$ws.Range("E37").Value = "'  -0.04%  "
$ws.Range("D38").Value = "'11.17"
$ws.Range("E38").Value = "'  +13.56%  "
$ws.Range("D39").Value = "'0.2759"
$ws.Range("E39").Value = "'  +3.75%  "
$ws.Range("D40").Value = "'0.09242"
$ws.Range("E40").Value = "'  +1.12%  "
$ws.Range("E41").Value = "'  +2.80%  "
$ws.Range("D42").Value = "'13.69"
$ws.Range("E42").Value = "'  +5.88%  "
$ws.Range("D43").Value = "'1.444"
$ws.Range("E43").Value = "'  +0.16%  "
$ws.Range("D44").Value = "'16.60"
$ws.Range("E44").Value = "'  +4.42%  "
$ws.Range("D45").Value = "'0.7077"
$ws.Range("E45").Value = "'  +3.47%  "
$ws.Range("D46").Value = "'2.540"
$ws.Range("E46").Value = "'  +3.99%  "
$ws.Range("D47").Value = "'4.130"
$ws.Range("E47").Value = "'  +1.28%  "
$ws.Range("D48").Value = "'1.000"
$ws.Range("E48").Value = "'  +0.05%  "
$ws.Range("D49").Value = "'0.08425"
$ws.Range("E49").Value = "'  +0.98%  "
$ws.Range("D50").Value = "'136.31"
$ws.Range("E50").Value = "'  +3.27%  "
$ws.Range("D51").Value = "'1.262"
$ws.Range("E51").Value = "'  +0.99%  "
